$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '69.643.08'
$ws.Range('E2').Value = '  +1.48%  '

$ws.Range('D3').Value = '3.890.97'
$ws.Range('E3').Value = '  +1.16%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '604.01'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.60%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '170.93'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.45%  '

$ws.Range('D7').Value = '3.891.11'
$ws.Range('E7').Value = '  +1.17%  '

$ws.Range('E8').Value = '  +0.06%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.536'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.21%  '

$ws.Range('E10').Value = '  +1.34%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.40'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.92%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000256'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +4.89%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '38.29'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.72%  '

$ws.Range('D15').Value = '4.544.65'
$ws.Range('E15').Value = '  +1.22%  '

$ws.Range('D16').Value = '3.883.07'
$ws.Range('E16').Value = '  +0.48%  '

$ws.Range('D17').Value = '69.671.56'
$ws.Range('E17').Value = '  +1.27%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '18.70'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +9.31%  '

$ws.Range('E19').Value = '  +0.95%  '

$ws.Range('E20').Value = '  -0.79%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.08'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.63%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '491.06'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.22%  '

$ws.Range('E23').Value = '  +3.77%  '

$ws.Range('E24').Value = '  +2.21%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '85.26'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.41%  '

$ws.Range('E26').Value = '  +2.68%  '

$ws.Range('E27').Value = '  +1.64%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.15'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +1.66%  '

$ws.Range('E29').Value = '  +0.21%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.99'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.30%  '

$ws.Range('D31').Value = '4.040.92'
$ws.Range('E31').Value = '  +1.11%  '

$ws.Range('E32').Value = '  +1.25%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.87'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.13%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '31.95'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.21%  '

$ws.Range('D35').Value = '3.858.18'
$ws.Range('E35').Value = '  +1.75%  '

$ws.Range('E36').Value = '  -0.04%  '

$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.41'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +14.71%  '

$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.11'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +3.92%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.143'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +2.01%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.04'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.60%  '

$ws.Range('E41').Value = '  +0.03%  '

$ws.Range('E42').Value = '  +2.93%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.09'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +5.40%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '433.95'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.41%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '48.08'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.90%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '8.69'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.51%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.000276'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +21.67%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0367'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.92%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '143.50'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.37%  '

$ws.Range('E51').Value = '  +4.22%  '
